$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 88; this shifts the existing rows 88-147 down to 89-148
# (and Excel auto-extends the used range from A1:R147 to A1:R148).
$ws.Rows.Item(88).Insert()

# Populate the newly inserted row 88 with the new "Alcachofa" price entry.
# (Same Mercado/Region/Categoria/Variedad/Calidad/Volumen/Unidad/Origen/Kg
# as the row that used to sit at 88, just a new date + new price figures.)
$ws.Range("A88").Value = 5
$ws.Range("B88").Value = "Macroferia Regional de Talca"
$ws.Range("C88").Value = "Maule"
$ws.Range("D88").Value = 45176
$ws.Range("E88").Value = 7
$ws.Range("F88").Value = 100112013
$ws.Range("G88").Value = "Alcachofa"
$ws.Range("H88").Value = "Madrigal"
$ws.Range("I88").Value = "Primera"
$ws.Range("J88").Value = 300
$ws.Range("K88").Value = 10000
$ws.Range("L88").Value = 10000
$ws.Range("M88").Value = 10000
$ws.Range("N88").Value = "$/caja 40 unidades"
$ws.Range("O88").Value = "Provincia del Elqu" + [char]0x00ED
$ws.Range("P88").Value = 250
$ws.Range("Q88").Value = 40
$ws.Range("R88").Value = "Hortaliza"
